$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the new "Logicalis Spain" role (two paragraphs: a Heading2 title
#    line with dates, and a Body Text description) immediately before the
#    existing "Infortec Consultores for Kyndryl" entry.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Infortec Consultores")) {
        $target = $p
        break
    }
}

$insertPoint = $d.Range($target.Range.Start, $target.Range.Start)

$newBlockXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
'<w:p><w:pPr><w:pStyle w:val="Ttulo2"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' + `
'<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Logicalis Spain</w:t></w:r>' + `
'<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:tab/><w:t>2022 – present</w:t></w:r>' + `
'</w:p>' + `
'<w:p><w:pPr><w:pStyle w:val="Textoindependiente"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' + `
'<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Cloud Administrator. </w:t></w:r>' + `
'<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Management</w:t></w:r>' + `
'<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> of OpenShift, both on premises and Azure. Automation with Ansible</w:t></w:r>' + `
'<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> and Python</w:t></w:r>' + `
'<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>.</w:t></w:r>' + `
'</w:p>' + `
'<w:p/>' + `
'</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($newBlockXml)

# The trailing empty <w:p/> in the inserted package is needed so the real
# "Infortec Consultores..." paragraph is not merged with our new content;
# it leaves behind one stray empty ("Normal" style) paragraph that we now
# remove.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "`r" -and $p.Style.NameLocal -eq "Normal") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# 2. The "Infortec Consultores for Kyndryl" role is no longer current, so its
#    date changes from "2022 - present" to just "2022".
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Infortec Consultores for Kyndryl" + [char]9 + "2022 – present", $true, $false, $false, $false, $false, $true, 1, $false, "Infortec Consultores for Kyndryl" + [char]9 + "2022", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. The NTT Managed Services Iberia role now has an end date of 2022.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("NTT Managed Services Iberia" + [char]9 + "2021 –", $true, $false, $false, $false, $false, $true, 1, $false, "NTT Managed Services Iberia" + [char]9 + "2021 – 2022", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Trim the virtualization sentence later in the "Accenture" entry.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" and support of operating systems (Windows Server, Linux, AIX, HP-UX, Solaris) and virtualization environments (VMware vSphere, HP Integrity VM, Sun LDOMs) in both production and development environments.", $true, $false, $false, $false, $false, $true, 1, $false, " and support of operating systems (Windows Server, Linux, AIX, HP-UX, Solaris) and virtualization in both production and development environments.", 2) | Out-Null
